# Updated flow and operating expense calculation
#
# G8 ("butter" row's electricity consumption) changes from 41000 to 100.
# This ripples through J8, K5:K9 (shared formula), M5:M9, N5:N9 and J10
# automatically on recalculation.
#
# The "Weight" column header (I4) is renamed to "SnsitivityParamt".
#
# The active selection moves from M5 to J5.
#
# Iterative calculation is turned on with a tighter max-change threshold
# (mirrors the workbook's calcPr iterateDelta="1E-4").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enable iterative calculation with a smaller max change (iterateDelta).
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# Rename the "Weight" header to "SnsitivityParamt".
$ws.Range("I4").Value = "SnsitivityParamt"

# Core data edit: electricity consumption for the "butter" row drops to 100.
$ws.Range("G8").Value = 100

# Move the active selection from M5 to J5.
$ws.Range("J5").Select() | Out-Null
